$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 21:05"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 1419478
$ws.Range("C4").Value = 10842
$ws.Range("E4").Value = 1034704
$ws.Range("G4").Value = 910
$ws.Range("H4").Value = 84335

# Row 10: Francia -> Francia
$ws.Range("B10").Value = 178060
$ws.Range("E10").Value = 92313

# Row 11: Alemania -> Alemania
$ws.Range("F11").Value = 1465

# Row 15: India -> India
$ws.Range("B15").Value = 78055
$ws.Range("C15").Value = 3763
$ws.Range("D15").Value = 26400
$ws.Range("E15").Value = 49104

# Row 16: Peru -> Canada
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 72196
$ws.Range("C16").Value = 1039
$ws.Range("D16").Value = 34916
$ws.Range("E16").Value = 31979
$ws.Range("F16").Value = 502
$ws.Range("G16").Value = 132
$ws.Range("H16").Value = 5301

# Row 17: Canada -> Peru
$ws.Range("A17").Value = "Peru"
$ws.Range("B17").Value = 72059
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 23324
$ws.Range("E17").Value = 46678
$ws.Range("F17").Value = 797
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 2057

# Row 59: Barein -> Barein
$ws.Range("B59").Value = 5816
$ws.Range("C59").Value = 285
$ws.Range("E59").Value = 3611
$ws.Range("F59").Value = 6
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 10

# Row 61: Moldavia -> Ghana
$ws.Range("A61").Value = "Ghana"
$ws.Range("B61").Value = 5408
$ws.Range("C61").Value = 281
$ws.Range("D61").Value = 514
$ws.Range("E61").Value = 4870
$ws.Range("F61").Value = 5
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 24

# Row 62: Afganistan -> Moldavia
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 5406
$ws.Range("C62").Value = 252
$ws.Range("D62").Value = 2176
$ws.Range("E62").Value = 3045
$ws.Range("F62").Value = 251
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 185

# Row 63: Ghana -> Afganistan
$ws.Range("A63").Value = "Afganistan"
$ws.Range("B63").Value = 5226
$ws.Range("C63").Value = 263
$ws.Range("D63").Value = 648
$ws.Range("E63").Value = 4446
$ws.Range("F63").Value = 7
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 132

# Row 101: Letonia -> Maldivas
$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 955
$ws.Range("C101").Value = 51
$ws.Range("D101").Value = 29
$ws.Range("E101").Value = 922
$ws.Range("H101").Value = 4

# Row 102: Republica de Chipre -> Letonia
$ws.Range("A102").Value = "Letonia"
$ws.Range("B102").Value = 951
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 627
$ws.Range("E102").Value = 305
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 19

# Row 103: Maldivas -> Republica de Chipre
$ws.Range("A103").Value = "Republica de Chipre"
$ws.Range("B103").Value = 905
$ws.Range("C103").Value = 2
$ws.Range("D103").Value = 449
$ws.Range("E103").Value = 439
$ws.Range("F103").Value = 10
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 17

# Row 106: Libano -> Libano
$ws.Range("F106").Value = 4

# Row 110: Costa Rica -> Costa Rica
$ws.Range("B110").Value = 815
$ws.Range("C110").Value = 11
$ws.Range("D110").Value = 527
$ws.Range("E110").Value = 281

# Row 114: Paraguay -> Mali
$ws.Range("A114").Value = "Mali"
$ws.Range("B114").Value = 758
$ws.Range("C114").Value = 28
$ws.Range("D114").Value = 412
$ws.Range("E114").Value = 302
$ws.Range("G114").Value = 4
$ws.Range("H114").Value = 44

# Row 115: Kenia -> Paraguay
$ws.Range("A115").Value = "Paraguay"
$ws.Range("B115").Value = 740
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 182
$ws.Range("E115").Value = 547
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 11

# Row 116: Mali -> Kenia
$ws.Range("A116").Value = "Kenia"
$ws.Range("B116").Value = 737
$ws.Range("C116").Value = 22
$ws.Range("D116").Value = 281
$ws.Range("E116").Value = 416
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 4

# Row 145: Liberia -> Liberia
$ws.Range("B145").Value = 213
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 101
$ws.Range("E145").Value = 92

# Row 165: Monaco -> Monaco
$ws.Range("D165").Value = 87
$ws.Range("E165").Value = 5

# Row 176: Angola -> Angola
$ws.Range("D176").Value = 14
$ws.Range("E176").Value = 29

# Row 193: Nueva Caledonia -> Belice
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Row 194: Belice -> Nueva Caledonia
$ws.Range("A194").Value = "Nueva Caledonia"
$ws.Range("D194").Value = 18
$ws.Range("H194").Value = 0

